# reproduc pipeline image added
#
# Applies the changes described by the diff:
#  1. Reflow/resize the "Query, Analyze, and Develop Manuscript ..." textbox
#     (id 33) and retext "with" -> "within".
#  2. Nudge the "Seamless Collaboration" textbox (id 36).
#  3. Add five new textboxes (Push / Pull / Push / Pull / Reproducibility...)
#     which become shape ids 41-45 (matching the diff) by first soaking up
#     every currently-unused id slot below 41 with throwaway textboxes,
#     adding the five real ones (which then land on 41-45), and finally
#     deleting the throwaway shapes again.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1. TextBox 32 (id 33) - "Query, Analyze, and Develop Manuscript ..."
# ---------------------------------------------------------------------------
$shQuery = $s.Shapes.Item(15)
$shQuery.Left   = 768.2014770507812
$shQuery.Top    = 36.14370346069336
$shQuery.Width  = 201.33111572265625
$shQuery.Height = 72.70315551757812
$shQuery.TextFrame.TextRange.Text = "Query, Analyze, and Develop Manuscript within Quarto Document"

# ---------------------------------------------------------------------------
# 2. TextBox 35 (id 36) - "Seamless Collaboration"
# ---------------------------------------------------------------------------
$shSeamless = $s.Shapes.Item(18)
$shSeamless.Left = 567.3451538085938
$shSeamless.Top  = 277.30181884765625

# ---------------------------------------------------------------------------
# 3. Five new textboxes (ids 41-45 in the target deck)
# ---------------------------------------------------------------------------
# The shape-id allocator always hands out the lowest id that is not
# currently in use on the slide, so to land exactly on 41-45 we first have
# to occupy every free slot from 2 up to 40 with scratch shapes, add the
# five real ones, then remove the scratch shapes (removing them afterwards
# does not renumber the shapes that already exist).
$usedIds = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $usedIds += $s.Shapes.Item($i).Id
}

$scratch = @()
for ($candidate = 2; $candidate -lt 41; $candidate++) {
    if ($usedIds -notcontains $candidate) {
        $tmp = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
        $scratch += $tmp
    }
}

# --- TextBox 40 (id 41) - "Push" ---
$tb41 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$tb41.Left   = 589.533203125
$tb41.Top    = 189.71205139160156
$tb41.Width  = 53.44504165649414
$tb41.Height = 29.081260681152344
$tb41.Rotation = 339.84710693359375
$tb41.Fill.Visible = $false
$tb41.TextFrame.WordWrap = $true
$tb41.TextFrame.AutoSize = 1
$tb41.TextFrame.TextRange.Text = "Push"
$tb41.Name = "TextBox 40"

# --- TextBox 41 (id 42) - "Pull" ---
$tb42 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$tb42.Left   = 601.471923828125
$tb42.Top    = 216.2370147705078
$tb42.Width  = 53.56157684326172
$tb42.Height = 29.081260681152344
$tb42.Rotation = 339.84710693359375
$tb42.Fill.Visible = $false
$tb42.TextFrame.WordWrap = $true
$tb42.TextFrame.AutoSize = 1
$tb42.TextFrame.TextRange.Text = "Pull"
$tb42.Name = "TextBox 41"

# --- TextBox 42 (id 43) - "Push" ---
$tb43 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$tb43.Left   = 600.7601318359375
$tb43.Top    = 329.15118408203125
$tb43.Width  = 58.28645706176758
$tb43.Height = 29.081260681152344
$tb43.Rotation = 12.694916725158691
$tb43.Fill.Visible = $false
$tb43.TextFrame.WordWrap = $true
$tb43.TextFrame.AutoSize = 1
$tb43.TextFrame.TextRange.Text = "Push"
$tb43.Name = "TextBox 42"

# --- TextBox 43 (id 44) - "Pull" ---
$tb44 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$tb44.Left   = 597.5626831054688
$tb44.Top    = 354.0142822265625
$tb44.Width  = 53.56157684326172
$tb44.Height = 29.081260681152344
$tb44.Rotation = 14.339966773986816
$tb44.Fill.Visible = $false
$tb44.TextFrame.WordWrap = $true
$tb44.TextFrame.AutoSize = 1
$tb44.TextFrame.TextRange.Text = "Pull"
$tb44.Name = "TextBox 43"

# --- TextBox 44 (id 45) - "Reproducibility results in Happy Healthy Frogs" ---
$tb45 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$tb45.Left   = 7.6862993240356445
$tb45.Top    = 186.31275939941406
$tb45.Width  = 198.45071411132812
$tb45.Height = 50.892208099365234
$tb45.Fill.Visible = $false
$tb45.TextFrame.WordWrap = $true
$tb45.TextFrame.AutoSize = 1
$tb45.TextFrame.TextRange.Text = "Reproducibility results in Happy Healthy Frogs"
$tb45.Name = "TextBox 44"

foreach ($tmp in $scratch) {
    $tmp.Delete()
}
